$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.802.56"
$ws.Range("E2").Value = "  +2.58%  "
$ws.Range("D3").Value = "1.876.13"
$ws.Range("E3").Value = "  +2.23%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "325.04"
$ws.Range("E5").Value = "  +0.40%  "
$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").Value = "0.4616"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("E8").Value = "  -0.35%  "
$ws.Range("D9").Value = "0.07873"
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("D10").Value = "0.9887"
$ws.Range("E10").Value = "  +2.59%  "
$ws.Range("D11").Value = "21.87"
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "1.891.11"
$ws.Range("E12").Value = "  +1.75%  "
$ws.Range("D13").Value = "7.005"
$ws.Range("E13").Value = "  +1.15%  "
$ws.Range("D14").Value = "5.706"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").Value = "0.06978"
$ws.Range("E15").Value = "  +1.98%  "
$ws.Range("D16").Value = "88.44"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "1.004"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("E18").Value = "  +0.95%  "
$ws.Range("E19").Value = "  +0.25%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").Value = "28.813.78"
$ws.Range("E21").Value = "  +2.52%  "
$ws.Range("D22").Value = "5.282"
$ws.Range("E22").Value = "  -0.76%  "
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D24").Value = "2.102"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").Value = "2.133.39"
$ws.Range("E25").Value = "  +4.52%  "
$ws.Range("D26").Value = "152.98"
$ws.Range("E26").Value = "  -1.17%  "
$ws.Range("D27").Value = "19.26"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").Value = "5.831"
$ws.Range("E28").Value = "  +2.32%  "
$ws.Range("D29").Value = "1.994"
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("D30").Value = "118.97"
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").Value = "0.09335"
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("D32").Value = "0.9219"
$ws.Range("E32").Value = "  -1.53%  "
$ws.Range("D33").Value = "5.307"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("D35").Value = "3.321"
$ws.Range("E35").Value = "  +0.39%  "
$ws.Range("D36").Value = "0.05795"
$ws.Range("E36").Value = "  -1.40%  "
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").Value = "0.02069"
$ws.Range("E38").Value = "  -2.88%  "
$ws.Range("D39").Value = "7.658"
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("D40").Value = "0.5634"
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("D42").Value = "9.797"
$ws.Range("D43").Value = "0.07217"
$ws.Range("E43").Value = "  -0.63%  "
$ws.Range("D44").Value = "11.80"
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("D45").Value = "0.5302"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").Value = "2.145"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").Value = "1.121"
$ws.Range("E47").Value = "  -3.22%  "
$ws.Range("D48").Value = "1.840"
$ws.Range("E48").Value = "  +0.71%  "
$ws.Range("D49").Value = "113.25"
$ws.Range("E49").Value = "  +0.65%  "
$ws.Range("E50").Value = "  +3.64%  "
$ws.Range("D51").Value = "1.003"
$ws.Range("E51").Value = "  +0.18%  "
